$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Apply "Bad" style (red fill / red text) to the two mis-sized panel rows ---
$ws.Range("B16:J17").Style = "Bad"

# Note of the sizing mistake, attached to the "panel wide" row
$ws.Range("I16").Value = "I made a size error"

# --- New section: replacement order for the corrected panels ---
$ws.Range("B23").Value = "Replacement order for panels"

$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 2610
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = "panel wide"
$ws.Range("F24").Value = "Panel"

$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 2610
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = "panel narrow"
$ws.Range("F25").Value = "Panel"

# the 9.86" (narrow) replacement note was typed before the 10.98" (wide) one
$ws.Range("H25").Value = "7155 Cut to size 9.86 x 53 inches"
$ws.Range("H24").Value = "7155 Cut to size 10.98 x 53 inches"

# --- Page setup: fit to page, landscape, scale 61% ---
$ws.PageSetup.Orientation = 2
$ws.PageSetup.Zoom = 61
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# --- Final selection as left by the edit ---
[void]$ws.Range("H26").Select()
